# Scheduled runner update: refresh market-board price snapshots (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-class Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 699.5
$ws.Range("J33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("N33").Value = -2458
$ws.Range("H58").Value = 1494.25
$ws.Range("J58").Value = 2242.8572
$ws.Range("L58").Value = 6728.571599999999
$ws.Range("N58").Value = -7028.571599999999
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H112").Value = 2275.8333
$ws.Range("J112").Value = 2513.75
$ws.Range("L112").Value = 7541.25
$ws.Range("N112").Value = -9757.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 4750
$ws.Range("I10").Value = 4750
$ws.Range("K10").Value = 4750
$ws.Range("M10").Value = -4580
$ws.Range("H13").Value = 925
$ws.Range("I13").Value = 925
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 925
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -781
$ws.Range("N13").ClearContents()
$ws.Range("H34").Value = 2523
$ws.Range("I34").Value = 2523
$ws.Range("K34").Value = 2523
$ws.Range("M34").Value = -2252
$ws.Range("H88").Value = 1503
$ws.Range("J88").Value = 1499
$ws.Range("L88").Value = 1499
$ws.Range("N88").Value = -2311
$ws.Range("H91").Value = 1503
$ws.Range("J91").Value = 1499
$ws.Range("L91").Value = 1499
$ws.Range("N91").Value = -4307
$ws.Range("H110").Value = 507.55554
$ws.Range("I110").Value = 321.125
$ws.Range("K110").Value = 321.125
$ws.Range("M110").Value = 1723.875
$ws.Range("H111").Value = 30644
$ws.Range("J111").Value = 30644
$ws.Range("L111").Value = 30644
$ws.Range("N111").Value = -38824
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 279.5
$ws.Range("I19").Value = 123.09091
$ws.Range("K19").Value = 123.09091
$ws.Range("M19").Value = 46.90909000000001
$ws.Range("H24").Value = 279.5
$ws.Range("I24").Value = 123.09091
$ws.Range("K24").Value = 123.09091
$ws.Range("M24").Value = 46.90909000000001
$ws.Range("H31").Value = 7039.846
$ws.Range("J31").Value = 8527.105
$ws.Range("L31").Value = 8527.105
$ws.Range("N31").Value = -9117.105
$ws.Range("H34").Value = 7039.846
$ws.Range("J34").Value = 8527.105
$ws.Range("L34").Value = 8527.105
$ws.Range("N34").Value = -8931.105
$ws.Range("H132").Value = 2347.5557
$ws.Range("I132").Value = 2347.5557
$ws.Range("K132").Value = 7042.6671
$ws.Range("M132").Value = -4512.6671
$ws.Range("H134").Value = 1916.3334
$ws.Range("I134").Value = 1800.6
$ws.Range("K134").Value = 5401.799999999999
$ws.Range("M134").Value = -2866.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 112743.164
$ws.Range("J4").Value = 1384.7142
$ws.Range("L4").Value = 4154.142599999999
$ws.Range("N4").Value = -4378.142599999999
$ws.Range("H12").Value = 126
$ws.Range("I12").Value = 17.5
$ws.Range("J12").Value = 165.45454
$ws.Range("K12").Value = 52.5
$ws.Range("L12").Value = 496.36362
$ws.Range("M12").Value = 120.5
$ws.Range("N12").Value = -842.3636200000001
$ws.Range("H16").Value = 439.5
$ws.Range("J16").Value = 439.5
$ws.Range("L16").Value = 1318.5
$ws.Range("N16").Value = -1664.5
$ws.Range("H21").Value = 333.33334
$ws.Range("J21").Value = 250
$ws.Range("L21").Value = 750
$ws.Range("N21").Value = -1096
$ws.Range("H29").Value = 10
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H33").Value = 158.21428
$ws.Range("I33").Value = 120.7
$ws.Range("K33").Value = 724.2
$ws.Range("M33").Value = -441.2
$ws.Range("H86").Value = 518.25
$ws.Range("I86").Value = 489.8
$ws.Range("K86").Value = 1469.4
$ws.Range("M86").Value = -283.4000000000001
$ws.Range("H89").Value = 518.25
$ws.Range("I89").Value = 489.8
$ws.Range("K89").Value = 4408.2
$ws.Range("M89").Value = 1519.8
$ws.Range("H94").Value = 7048.8887
$ws.Range("I94").Value = 4024
$ws.Range("J94").Value = 7427
$ws.Range("K94").Value = 12072
$ws.Range("L94").Value = 22281
$ws.Range("M94").Value = -11396
$ws.Range("N94").Value = -23633

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1328.3846
$ws.Range("I102").Value = 1328.3846
$ws.Range("K102").Value = 1328.3846
$ws.Range("M102").Value = 293.6153999999999
$ws.Range("H112").Value = 20293
$ws.Range("J112").Value = 20293
$ws.Range("L112").Value = 20293
$ws.Range("N112").Value = -22509

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9980.333000000001
$ws.Range("I7").Value = 9962.666999999999
$ws.Range("J7").Value = 9998
$ws.Range("K7").Value = 9962.666999999999
$ws.Range("L7").Value = 9998
$ws.Range("M7").Value = -9850.666999999999
$ws.Range("N7").Value = -10222
$ws.Range("H12").Value = 701
$ws.Range("J12").Value = 701
$ws.Range("L12").Value = 701
$ws.Range("N12").Value = -1041
$ws.Range("H20").Value = 51649.95
$ws.Range("J20").Value = 53578.895
$ws.Range("L20").Value = 53578.895
$ws.Range("N20").Value = -54030.895
$ws.Range("H22").Value = 901.25
$ws.Range("I22").Value = 830.7143
$ws.Range("K22").Value = 830.7143
$ws.Range("M22").Value = -535.7143
$ws.Range("H27").Value = 901.25
$ws.Range("I27").Value = 830.7143
$ws.Range("K27").Value = 830.7143
$ws.Range("M27").Value = -723.7143
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H68").Value = 8317.454
$ws.Range("I68").Value = 7499.8
$ws.Range("K68").Value = 7499.8
$ws.Range("M68").Value = -6750.8
$ws.Range("H71").Value = 8317.454
$ws.Range("I71").Value = 7499.8
$ws.Range("K71").Value = 37499
$ws.Range("M71").Value = -33755
$ws.Range("H100").Value = 8333.111000000001
$ws.Range("J100").Value = 9124.875
$ws.Range("L100").Value = 9124.875
$ws.Range("N100").Value = -10206.875
$ws.Range("H126").Value = 9980.333000000001
$ws.Range("I126").Value = 9962.666999999999
$ws.Range("J126").Value = 9998
$ws.Range("K126").Value = 29888.001
$ws.Range("L126").Value = 29994
$ws.Range("M126").Value = -27418.001
$ws.Range("N126").Value = -34934
$ws.Range("H139").Value = 250000
$ws.Range("J139").Value = 250000
$ws.Range("L139").Value = 250000
$ws.Range("N139").Value = -260280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 5001.3335
$ws.Range("I7").Value = 5001.3335
$ws.Range("K7").Value = 5001.3335
$ws.Range("M7").Value = -4888.3335
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H104").Value = 26075.428
$ws.Range("J104").Value = 26075.428
$ws.Range("L104").Value = 26075.428
$ws.Range("N104").Value = -33063.428
